$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old orphan "Systolic" label (row 7) so it drops out of use ---
$ws.Range("C7").Value2 = ""

# --- Establish new shared-string creation order: 2-cores, single-cores, total, us ---
$ws.Range("C9").Value2  = "2-cores"
$ws.Range("C3").Value2  = "single-cores"
$ws.Range("I4").Value2  = "total"
$ws.Range("C4").Value2  = "us"

# --- Row 4: header row for the "single-cores" table ---
$ws.Range("D4").Value2 = "CONV1"
$ws.Range("E4").Value2 = "CONV2"
$ws.Range("F4").Value2 = "FC1"
$ws.Range("G4").Value2 = "FC2"
$ws.Range("H4").Value2 = "FC3"

# --- Row 5: Base values (single-cores) ---
$ws.Range("C5").Value2 = "Base"
$ws.Range("D5").Value2 = 1860292
$ws.Range("E5").Value2 = 1300248
$ws.Range("F5").Value2 = 215341
$ws.Range("G5").Value2 = 45394
$ws.Range("H5").Value2 = 4761
$ws.Range("I5").Value2 = 3424721

# --- Row 6: Row based values (single-cores) ---
$ws.Range("C6").Value2 = "Row based"
$ws.Range("D6").Value2 = 542124
$ws.Range("E6").Value2 = 59945
$ws.Range("F6").Value2 = 121786
$ws.Range("G6").Value2 = 26239
$ws.Range("H6").Value2 = 3313
$ws.Range("I6").Value2 = 262148

# --- Row 10: header row for the "2-cores" table ---
$ws.Range("C10").Value2 = "us"
$ws.Range("D10").Value2 = "CONV1"
$ws.Range("E10").Value2 = "CONV2"
$ws.Range("F10").Value2 = "FC1"
$ws.Range("G10").Value2 = "FC2"
$ws.Range("H10").Value2 = "FC3"
$ws.Range("I10").Value2 = "total"

# --- Row 11: Base values (2-cores) ---
$ws.Range("C11").Value2 = "Base"
$ws.Range("D11").Value2 = 1461003
$ws.Range("E11").Value2 = 1010759
$ws.Range("F11").Value2 = 180305
$ws.Range("G11").Value2 = 38025
$ws.Range("H11").Value2 = 3941
$ws.Range("I11").Value2 = 2693829

# --- Row 12: Row based values (2-cores) ---
$ws.Range("C12").Value2 = "Row based"
$ws.Range("D12").Value2 = 42052
$ws.Range("E12").Value2 = 46022
$ws.Range("F12").Value2 = 97821
$ws.Range("G12").Value2 = 20799
$ws.Range("H12").Value2 = 2808
$ws.Range("I12").Value2 = 209361

# --- Update the active selection to match target ---
[void]$ws.Range("H7").Select()
